$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Sex" column's code changed from SEX to MONSEX.
$ws.Range("D1").Value2 = "Sex (MONSEX)"

# A new "Total maximum statuary prison requirement (STAMAX)" column is being
# added right after the existing "Total minimum statuary prison requirement
# (STAMIN)" column (G), pushing everything from H onward one column to the
# right.
$ws.Columns("H").Insert()

# Fill in the header for the newly inserted column.
$ws.Range("H1").Value2 = "Total maximum statuary prison requirement (STAMAX)"

# Leave the active selection on the newly added header cell, matching the
# author's last edit location.
$ws.Range("H1").Select()
